$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width changed from 15.42578125 to 15.7109375 (stored units).
# The COM ColumnWidth setter here quantizes to 1/6-character steps, so we
# pick the input that lands closest to the target stored width.
$ws.Columns.Item(1).ColumnWidth = 14.8

# Update cell values (A1:B32)
$ws.Range("A1").Value = -0.3879292809067465
$ws.Range("B1").Value = 0.38669382708897615
$ws.Range("A2").Value = -0.26219023848929801
$ws.Range("B2").Value = 0.25917392305462172
$ws.Range("A3").Value = -0.15622098551378727
$ws.Range("B3").Value = 0.15532093016392778
$ws.Range("A4").Value = -0.14332093026568415
$ws.Range("B4").Value = 0.14251707231617772
$ws.Range("A5").Value = -0.13651707267725577
$ws.Range("B5").Value = 0.13491051500186124
$ws.Range("A6").Value = -0.033942045191999171
$ws.Range("B6").Value = 0.033927206050355263
$ws.Range("A7").Value = -0.013927206495335298
$ws.Range("B7").Value = 0.013917951223296754
$ws.Range("A8").Value = 0.0060820483310868312
$ws.Range("B8").Value = -0.006098392873437497
$ws.Range("A9").Value = 0.01209839250156719
$ws.Range("B9").Value = -0.012126874499260865
$ws.Range("A10").Value = -0.040057693823037255
$ws.Range("B10").Value = 0.040032782823743673
$ws.Range("A11").Value = -0.03553278318660702
$ws.Range("B11").Value = 0.03549545661684661
$ws.Range("A12").Value = -0.029495456988534396
$ws.Range("B12").Value = 0.029390127054937487
$ws.Range("A13").Value = -0.023390127430246821
$ws.Range("B13").Value = 0.023366289474902047
$ws.Range("A14").Value = -0.011366289882548841
$ws.Range("B14").Value = 0.01136139087411614
$ws.Range("A15").Value = -0.005361391250897185
$ws.Range("B15").Value = 0.0053602187825241288
$ws.Range("A16").Value = 0.00063978084046789618
$ws.Range("B16").Value = -0.00064017554585849012
$ws.Range("A17").Value = 0.0066401751689584287
$ws.Range("B17").Value = -0.0066422995748594005
$ws.Range("A18").Value = -0.11450358896206936
$ws.Range("B18").Value = 0.11430523021353167
$ws.Range("A19").Value = -0.10530523057194019
$ws.Range("B19").Value = 0.10370604253052518
$ws.Range("A20").Value = -0.018014201995923074
$ws.Range("B20").Value = 0.018004436896695175
$ws.Range("A21").Value = -0.0090044372692457131
$ws.Range("B21").Value = 0.0089999996270750948
$ws.Range("A22").Value = -0.093952937908863277
$ws.Range("B22").Value = 0.093638118694084227
$ws.Range("A23").Value = -0.084638119065450823
$ws.Range("B23").Value = 0.08412757735491283
$ws.Range("A24").Value = -0.042127577904695457
$ws.Range("B24").Value = 0.041999999447278391
$ws.Range("A25").Value = -0.09496847018120036
$ws.Range("B25").Value = 0.094720510246489198
$ws.Range("A26").Value = -0.088720510619221926
$ws.Range("B26").Value = 0.088400919940955447
$ws.Range("A27").Value = -0.08240092031567281
$ws.Range("B27").Value = 0.081307647578885867
$ws.Range("A28").Value = -0.075307647961865065
$ws.Range("B28").Value = 0.074550096590702886
$ws.Range("A29").Value = -0.062550097010452888
$ws.Range("B29").Value = 0.062174436722298054
$ws.Range("A30").Value = -0.042174437186568436
$ws.Range("B30").Value = 0.042020423956663944
$ws.Range("A31").Value = -0.027020424398731535
$ws.Range("B31").Value = 0.027000831906768497
$ws.Range("A32").Value = -0.0060008323809155328
$ws.Range("B32").Value = 0.0059999996041142367
